$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the PREPROD claim numbers in column B (rows 4-7)
# Use a leading apostrophe (quote prefix) so Excel keeps these as text,
# preserving leading zeros and trailing spaces, matching the existing
# quotePrefix cell style already applied to this range.
$ws.Range("B4").Value = "'0420172008381   "
$ws.Range("B5").Value = "'0420172008380"
$ws.Range("B6").Value = "'1220170301405"
$ws.Range("B7").Value = "'1120170200926"

# Delete row 8 which previously held the 5th PREPROD claim entry
$ws.Rows.Item(8).Delete()

# Move the active selection like the author's final recorded state
$ws.Range("I4").Select()
